$d = $word.ActiveDocument

# 1. Replace the title text
$d.Content.Find.Execute(
    "Section 41 –  Deploying A Vue JS Application",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Section 21 – Roundup And Next Steps", 2)

# 2. Remove the "Chapter 635 –  Chapter 637" paragraph and the following
#    empty paragraph (the one that only carries a tab stop), merging the
#    remaining content back together.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Chapter\s*635") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $startRange = $target.Range.Start
    # the paragraph immediately after the "Chapter 635..." one is the
    # empty paragraph that only holds a tab-stop definition
    $nextPara = $target.Next()
    $endRange = $nextPara.Range.End

    $r = $d.Range($startRange, $endRange)
    $r.Delete()
}
